# Fix Training Data Issue (#48)
#
# The "Date" column (BF) on this team game-log sheet was populated with the
# wrong text because of how NBA.com showed stats spanning two calendar
# years (season "2007-08") -- the literal string "5-29-2007-08" got stored
# instead of the actual game date. Correct every data row (row 2 through
# the last row, 31) to the real ISO-style date text "2008-05-29".
#
# NOTE: "2008-05-29" looks like a date to Excel's automatic type detection,
# so it is written with a leading apostrophe to force it to stay a literal
# text value (matching the original cell's text/string type) instead of
# silently turning into a date serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "5-29-2007-08"
$newText = "2008-05-29"

$dateRange = $ws.Range("BF2:BF31")

foreach ($cell in $dateRange.Cells) {
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = "'" + $newText
    }
}
